$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Add "under com\example\BookshelfTop" as a new run right after the
#    "Javadoc files are located in the Javadoc folder " run, and move
#    the hidden "_GoBack" bookmark to sit right after the new text.
# ------------------------------------------------------------------

$javadocPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Javadoc files are located in the Javadoc folder `r") {
        $javadocPara = $p
        break
    }
}

# Remove the old hidden "_GoBack" bookmark first (it currently sits
# after the "Click "Bookshelf"" run further down the document) so its
# id is free and Word's "last edit" bookmark can be re-created in the
# new spot below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

if ($javadocPara -ne $null) {
    $jRng = $javadocPara.Range

    $fragJavadoc = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' +
        '<w:p w14:paraId="2728596E" w14:textId="3B6FA7C5" w:rsidR="00C17981" w:rsidRDefault="00C17981" w:rsidP="00C17981" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
        '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="24"/></w:numPr></w:pPr>' +
        '<w:r><w:t xml:space="preserve">Javadoc files are located in the Javadoc folder </w:t></w:r>' +
        '<w:r><w:t>under com\example\BookshelfTop</w:t></w:r>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
        '<w:bookmarkEnd w:id="0"/>' +
        '</w:p>' +
        '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $jRng.InsertXML($fragJavadoc)
}

# ------------------------------------------------------------------
# 2) Shrink the VML screenshot ("Capture") from 249.3pt to 249pt wide.
# ------------------------------------------------------------------

$picturePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.WordOpenXML.IndexOf("width:249.3pt;height:201.75pt") -ge 0) {
        $picturePara = $p
        break
    }
}

if ($picturePara -ne $null) {
    $pRng = $picturePara.Range
    $xml = $pRng.WordOpenXML
    $xml2 = $xml.Replace("width:249.3pt;height:201.75pt", "width:249pt;height:201.75pt")
    if ($xml2 -ne $xml) {
        $pRng.InsertXML($xml2)
    }
}
